# Add the "Ophold" (residence permit) mapping worksheet after "Lande"
$wb = $excel.ActiveWorkbook
$landeSheet = $wb.Worksheets.Item("Lande")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $landeSheet)
$ws.Name = "Ophold"

# Column A ("Opholdsgrundlag") values, row 1..25, written first so the
# shared-string pool fills in the same order as the source workbook.
$opholdsgrundlag = @(
    'Opholdsgrundlag',
    'Familiesammenføring, Andre familiemedlemmer, Refererer til flygtning',
    'Ukraine (særlov)',
    'Familiesammenføring, Andre familiemedlemmer, Refererer til dansk/nordisk person',
    'Familiesammenføring, Mindreårige børn, Refererer til flygtning',
    'Familiesammenføring, Andre familiemedlemmer, Uoplyst referenceperson',
    'Det øvrige opholdsområde, Adoption',
    'Asyl, Andet grundlag',
    'Familiesammenføring, Ægteskab eller fast samlivsforhold, Refererer til udlænding, men ikke flygtning',
    'Familiesammenføring, Ægteskab eller fast samlivsforhold, Refererer til flygtning',
    'EU/EØS, Øvrige grunde',
    'EU/EØS, Uddannelse',
    'EU/EØS, Lønarbejde',
    'Familiesammenføring, Mindreårige børn, Refererer til andre end flygtning',
    'Asyl, Flygtningestatus',
    'Familiesammenføring, Mindreårige børn, Uoplyst referenceperson',
    'Studie mv., Praktikanter',
    'Studie mv., Au pair',
    'Familiesammenføring, Ægteskab eller fast samlivsforhold, Refererer til dansk/nordisk person',
    'Det øvrige opholdsområde, Øvrige grunde',
    'Familiesammenføring, Ægteskab eller fast samlivsforhold, Uoplyst referenceperson',
    'Studie mv., Øvrige grunde',
    'EU/EØS, Familiemedlemmer',
    'Studie mv., Uddannelse',
    'Erhverv'
)

# Column B ("Opholdstype") values, row 1..25.
$opholdstype = @(
    'Opholdstype',
    'Familiesammenføring',
    'Flygtningestatus/asyl',
    'Familiesammenføring',
    'Familiesammenføring',
    'Familiesammenføring',
    'Adoption',
    'Flygtningestatus/asyl',
    'Familiesammenføring',
    'Familiesammenføring',
    'Øvrige grunde',
    'Uddannelse/praktik',
    'Lønarbejde/erhverv',
    'Familiesammenføring',
    'Flygtningestatus/asyl',
    'Familiesammenføring',
    'Uddannelse/praktik',
    'Uddannelse/praktik',
    'Familiesammenføring',
    'Øvrige grunde',
    'Familiesammenføring',
    'Uddannelse/praktik',
    'Familiesammenføring',
    'Uddannelse/praktik',
    'Lønarbejde/erhverv'
)

for ($i = 0; $i -lt $opholdsgrundlag.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $opholdsgrundlag[$i]
}
for ($i = 0; $i -lt $opholdstype.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $opholdstype[$i]
}

# Column widths (approximate Excel auto-fit sizing used in the source file).
$ws.Columns.Item(1).ColumnWidth = 82.21875
$ws.Columns.Item(2).ColumnWidth = 19.5546875

# Turn the range into a proper table ("Ophold"), styled like TableStyleMedium6.
$tableRange = $ws.Range("A1:B25")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Ophold"
$lo.TableStyle = "TableStyleMedium6"

# Leave selection on B15, matching the source file state.
$ws.Range("B15").Select()
